$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3865.6667
$ws.Range("I76").Value = 3424.75
$ws.Range("J76").Value = 4747.5
$ws.Range("K76").Value = 3424.75
$ws.Range("L76").Value = 4747.5
$ws.Range("M76").Value = -3109.75
$ws.Range("N76").Value = -5377.5

$ws.Range("H79").Value = 3865.6667
$ws.Range("I79").Value = 3424.75
$ws.Range("J79").Value = 4747.5
$ws.Range("K79").Value = 3424.75
$ws.Range("L79").Value = 4747.5
$ws.Range("M79").Value = -2332.75
$ws.Range("N79").Value = -6931.5

$ws.Range("H80").Value = 716
$ws.Range("I80").Value = 282.33334
$ws.Range("J80").Value = 1236.4
$ws.Range("K80").Value = 847.0000200000001
$ws.Range("L80").Value = 3709.2
$ws.Range("M80").Value = 150.9999799999999
$ws.Range("N80").Value = -5705.200000000001

$ws.Range("H83").Value = 716
$ws.Range("I83").Value = 282.33334
$ws.Range("J83").Value = 1236.4
$ws.Range("K83").Value = 2541.00006
$ws.Range("L83").Value = 11127.6
$ws.Range("M83").Value = 2450.99994
$ws.Range("N83").Value = -21111.6

$ws.Range("H88").Value = 1318.125
$ws.Range("I88").Value = 886.5
$ws.Range("J88").Value = 1749.75
$ws.Range("K88").Value = 886.5
$ws.Range("L88").Value = 1749.75
$ws.Range("M88").Value = -480.5
$ws.Range("N88").Value = -2561.75

$ws.Range("H91").Value = 1318.125
$ws.Range("I91").Value = 886.5
$ws.Range("J91").Value = 1749.75
$ws.Range("K91").Value = 886.5
$ws.Range("L91").Value = 1749.75
$ws.Range("M91").Value = 517.5
$ws.Range("N91").Value = -4557.75

$ws.Range("H96").Value = 625

$ws.Range("H103").Value = 831.6667
$ws.Range("I103").Value = 747.5
$ws.Range("K103").Value = 2242.5
$ws.Range("M103").Value = -1656.5

$ws.Range("H129").Value = 3753.4443
$ws.Range("I129").Value = 1533.3334
$ws.Range("K129").Value = 4600.0002
$ws.Range("M129").Value = 399.9997999999996

$ws.Range("H132").Value = 1246.9546
$ws.Range("I132").Value = 908.82355
$ws.Range("K132").Value = 2726.47065
$ws.Range("M132").Value = -196.4706499999998

$ws.Range("H135").Value = 1326.5714
$ws.Range("I135").Value = 1172.091
$ws.Range("J135").Value = 1893
$ws.Range("K135").Value = 10548.819
$ws.Range("L135").Value = 17037
$ws.Range("M135").Value = -8013.819
$ws.Range("N135").Value = -22107

$ws.Range("H138").Value = 1663.9565
$ws.Range("I138").Value = 454.2
$ws.Range("K138").Value = 1362.6
$ws.Range("M138").Value = 3777.4

$ws.Range("H141").Value = 4695.0557
$ws.Range("I141").Value = 4912.706
$ws.Range("J141").Value = 995
$ws.Range("K141").Value = 14738.118
$ws.Range("L141").Value = 2985
$ws.Range("M141").Value = -9558.118
$ws.Range("N141").Value = -13345

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 60000
$ws.Range("J92").Value = 60000
$ws.Range("L92").Value = 60000
$ws.Range("N92").Value = -64992

$ws.Range("H102").Value = 1986.3334
$ws.Range("I102").Value = 1986.3334
$ws.Range("K102").Value = 1986.3334
$ws.Range("M102").Value = -364.3334

$ws.Range("H135").Value = 48357.25
$ws.Range("J135").Value = 48357.25
$ws.Range("L135").Value = 48357.25
$ws.Range("N135").Value = -58497.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 24432.572
$ws.Range("J75").Value = 67499.5
$ws.Range("L75").Value = 67499.5
$ws.Range("N75").Value = -69371.5

$ws.Range("H78").Value = 24432.572
$ws.Range("J78").Value = 67499.5
$ws.Range("L78").Value = 202498.5
$ws.Range("N78").Value = -211858.5

$ws.Range("H82").Value = 29555
$ws.Range("I82").Value = 8665
$ws.Range("K82").Value = 8665
$ws.Range("M82").Value = -8282

$ws.Range("H85").Value = 29555
$ws.Range("I85").Value = 8665
$ws.Range("K85").Value = 8665
$ws.Range("M85").Value = -7339

$ws.Range("H86").Value = 1848.2858
$ws.Range("I86").Value = 2022.25
$ws.Range("J86").Value = 1616.3334
$ws.Range("K86").Value = 2022.25
$ws.Range("L86").Value = 1616.3334
$ws.Range("M86").Value = -899.25
$ws.Range("N86").Value = -3862.3334

$ws.Range("H89").Value = 1848.2858
$ws.Range("I89").Value = 2022.25
$ws.Range("J89").Value = 1616.3334
$ws.Range("K89").Value = 10111.25
$ws.Range("L89").Value = 8081.666999999999
$ws.Range("M89").Value = -4495.25
$ws.Range("N89").Value = -19313.667

$ws.Range("H94").Value = 1175.125
$ws.Range("I94").Value = 1215.7778
$ws.Range("K94").Value = 1215.7778
$ws.Range("M94").Value = -764.7778000000001

$ws.Range("H99").Value = 2673.4443
$ws.Range("I99").Value = 2514.6
$ws.Range("K99").Value = 2514.6
$ws.Range("M99").Value = -1016.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 36880
$ws.Range("J87").Value = 36880
$ws.Range("L87").Value = 36880
$ws.Range("N87").Value = -39252

$ws.Range("H90").Value = 36880
$ws.Range("J90").Value = 36880
$ws.Range("L90").Value = 110640
$ws.Range("N90").Value = -122496

$ws.Range("H134").Value = 2760.8
$ws.Range("I134").Value = 1349.7273
$ws.Range("K134").Value = 4049.1819
$ws.Range("M134").Value = -1514.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 502999.25
$ws.Range("J11").Value = 3999
$ws.Range("L11").Value = 11997
$ws.Range("N11").Value = -12277

$ws.Range("H12").Value = 464.16666
$ws.Range("I12").Value = 443.5
$ws.Range("K12").Value = 1330.5
$ws.Range("M12").Value = -1157.5

$ws.Range("H34").Value = 2340.25
$ws.Range("J34").Value = 2878.6667
$ws.Range("L34").Value = 8636.000100000001
$ws.Range("N34").Value = -8804.000100000001

$ws.Range("H39").Value = 3062.5
$ws.Range("J39").Value = 3666.6667
$ws.Range("L39").Value = 11000.0001
$ws.Range("N39").Value = -11588.0001

$ws.Range("H55").Value = 1558.4
$ws.Range("J55").Value = 1720.5555
$ws.Range("L55").Value = 5161.666499999999
$ws.Range("N55").Value = -5515.666499999999

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("N106").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3623.5454
$ws.Range("I80").Value = 2538.75
$ws.Range("J80").Value = 4243.4287
$ws.Range("K80").Value = 2538.75
$ws.Range("L80").Value = 4243.4287
$ws.Range("M80").Value = -1540.75
$ws.Range("N80").Value = -6239.4287

$ws.Range("H83").Value = 3623.5454
$ws.Range("I83").Value = 2538.75
$ws.Range("J83").Value = 4243.4287
$ws.Range("K83").Value = 12693.75
$ws.Range("L83").Value = 21217.1435
$ws.Range("M83").Value = -7701.75
$ws.Range("N83").Value = -31201.1435

$ws.Range("H102").Value = 541.1
$ws.Range("I102").Value = 516.9474
$ws.Range("K102").Value = 516.9474
$ws.Range("M102").Value = 1105.0526

$ws.Range("H132").Value = 2614.8333
$ws.Range("I132").Value = 2537.8
$ws.Range("K132").Value = 7613.400000000001
$ws.Range("M132").Value = -5083.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1436.6
$ws.Range("I93").Value = 1570.875
$ws.Range("J93").Value = 899.5
$ws.Range("K93").Value = 1570.875
$ws.Range("L93").Value = 899.5
$ws.Range("M93").Value = -322.875
$ws.Range("N93").Value = -3395.5

$ws.Range("H98").Value = 30000
$ws.Range("I98").Value = 30000
$ws.Range("K98").Value = 30000
$ws.Range("M98").Value = -27005

$ws.Range("H100").Value = 300.5
$ws.Range("I100").Value = 300.5
$ws.Range("K100").Value = 300.5
$ws.Range("M100").Value = 240.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 20287.143
$ws.Range("J101").Value = 20287.143
$ws.Range("L101").Value = 20287.143
$ws.Range("N101").Value = -26777.143

$ws.Range("H104").Value = 20499.5
$ws.Range("J104").Value = 20499.5
$ws.Range("L104").Value = 20499.5
$ws.Range("N104").Value = -27487.5

$ws.Range("H132").Value = 1128.1818
$ws.Range("I132").Value = 1145.5555
$ws.Range("K132").Value = 3436.6665
$ws.Range("M132").Value = -906.6664999999998
